$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dish = "[teste] Frango ao molho de queijos"

$data = @(
    @($dish, "File de peito de frango", 300, "g"),
    @($dish, "Arroz Parboilizado", 150, "Un"),
    @($dish, "Queijo Mussarela fatiado", 50, "Kg"),
    @($dish, "Queijo Parmesao", 50, "g")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
